$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.156.41'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.866.83'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.41%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.82'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4705'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3916'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.91'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07979'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9967'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.69'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.981'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.269'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.842.57'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.41'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.98%  '

$ws.Range("E18").Value = '  +0.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06603'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.75'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.72%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9994'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.147.91'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.436'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.90%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.03'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.285'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.13%  '

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.068.29'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.07'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.77'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.128'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.498'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.75'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9736'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09492'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.94%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.571'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.52%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.374'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.62%  '

$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.339'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.74%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02271'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.91%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06107'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.402'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.174'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5973'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.45%  '

$ws.Range("E42").Value = '  -0.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1882'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.33'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.276'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5618'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.15'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.980'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.70%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06860'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.71'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.64%  '

$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.054'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.25%  '
